# edit.ps1 -- applies the "Fri May  3 20:16:37 UTC 2024" cryptos-list refresh
# (Coin / Link / Price / Volume(1h) table on Sheet1, rows 2-51).
#
# The Price column (D) holds values that look numeric (e.g. "62.261.76",
# "0.467") but must stay TEXT cells, matching the original workbook. Excel
# auto-converts a numeric-looking string assigned to Range.Value into a real
# number, so for every Price cell we touch we first force
# NumberFormat = "@" (Text) on that cell -- exactly like pre-formatting a
# cell as Text before typing into it in the Excel UI -- which keeps the
# write as a text value instead of a number.
#
# Rows 28 and 30 also swapped which coin they list (RenderToken now sits at
# row 28, FirstDigitalUSD at row 30), so those two rows rewrite Coin (B) and
# Link (C) as well as Price/Volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.261.76'
$ws.Range("E2").Value = '  +4.66%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.080.62'
$ws.Range("E3").Value = '  +2.53%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.62'
$ws.Range("E5").Value = '  +3.54%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.69'
$ws.Range("E6").Value = '  +3.76%  '

# Row 7
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.071.24'
$ws.Range("E8").Value = '  +2.67%  '

# Row 9
$ws.Range("E9").Value = '  +1.14%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.141'
$ws.Range("E10").Value = '  +5.97%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.76'
$ws.Range("E11").Value = '  +11.31%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.467'
$ws.Range("E12").Value = '  +2.15%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000242'
$ws.Range("E13").Value = '  +4.65%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.43'
$ws.Range("E14").Value = '  +4.81%  '

# Row 15
$ws.Range("E15").Value = '  +0.57%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.584.73'
$ws.Range("E16").Value = '  +2.54%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.30'
$ws.Range("E17").Value = '  -0.01%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.076.31'
$ws.Range("E18").Value = '  +2.58%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '62.171.45'
$ws.Range("E19").Value = '  +4.85%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '449.51'
$ws.Range("E20").Value = '  +4.05%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.04'
$ws.Range("E21").Value = '  +2.58%  '

# Row 22
$ws.Range("E22").Value = '  +1.81%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.47'
$ws.Range("E23").Value = '  +4.75%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.75'
$ws.Range("E24").Value = '  +2.99%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.86'
$ws.Range("E25").Value = '  +1.22%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.11%  '

# Row 27
$ws.Range("E27").Value = '  +4.04%  '

# Row 28
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.29'
$ws.Range("E28").Value = '  +5.46%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.68'
$ws.Range("E29").Value = '  +4.80%  '

# Row 30
$ws.Range("B30").Value = 'FirstDigitalUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.13%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.84'
$ws.Range("E31").Value = '  +11.79%  '

# Row 32
$ws.Range("E32").Value = '  +13.68%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.93'
$ws.Range("E33").Value = '  +4.43%  '

# Row 34
$ws.Range("E34").Value = '  +4.26%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0797'
$ws.Range("E35").Value = '  +4.64%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.06'
$ws.Range("E36").Value = '  +1.87%  '

# Row 37
$ws.Range("E37").Value = '  +5.52%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.43'
$ws.Range("E38").Value = '  +2.65%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.97'
$ws.Range("E39").Value = '  +9.12%  '

# Row 40
$ws.Range("E40").Value = '  +1.34%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '425.11'
$ws.Range("E41").Value = '  +5.88%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.909.74'
$ws.Range("E42").Value = '  +5.49%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0369'
$ws.Range("E43").Value = '  +4.53%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.278'
$ws.Range("E44").Value = '  +9.89%  '

# Row 45
$ws.Range("E45").Value = '  +0.76%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.16'
$ws.Range("E46").Value = '  +7.51%  '

# Row 47
$ws.Range("E47").Value = '  +0.05%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.47'
$ws.Range("E48").Value = '  +1.63%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.77'
$ws.Range("E49").Value = '  -2.98%  '

# Row 50
$ws.Range("E50").Value = '  +0.31%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.60'
$ws.Range("E51").Value = '  +4.68%  '
